$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 15
$ws.Range("B1").Value = 2.811738252639771
$ws.Range("C1").Value = 2.547181367874146
$ws.Range("D1").Value = 2.849354982376099
$ws.Range("E1").Value = 15
